$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "accountId"
$ws.Range("C1").Value = "region"

# Update data row 2
$ws.Range("A2").Value = "shubhamTest"
$ws.Range("B2").Value = 3438
$ws.Range("C2").Value = "Mumbai"

# Remove row 3 entirely (was testResource/ShubhamTest duplicate row)
$ws.Range("A3:C3").EntireRow.Delete()
